$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Change 1: "barkod_espp" -> split into two runs ("barkod_esp" + "p") with the
# existing _GoBack bookmark now sitting between them instead of after the
# single run.
# ---------------------------------------------------------------------------
$target1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*barkod_espp*") {
        $target1 = $p
        break
    }
}

if ($target1 -eq $null) {
    throw "Could not find paragraph containing 'barkod_espp'"
}

$xml1 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00B26294" w:rsidRPr="00203833" w:rsidRDefault="00B26294" w:rsidP="00C23685">' +
          '<w:pPr>' +
            '<w:rPr>' +
              '<w:highlight w:val="yellow"/>' +
              '<w:lang w:val="en-US"/>' +
            '</w:rPr>' +
          '</w:pPr>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r w:rsidRPr="00203833">' +
            '<w:rPr>' +
              '<w:highlight w:val="yellow"/>' +
              '<w:lang w:val="en-US"/>' +
            '</w:rPr>' +
            '<w:t>barkod_esp</w:t>' +
          '</w:r>' +
          '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
          '<w:bookmarkEnd w:id="0"/>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:highlight w:val="yellow"/>' +
              '<w:lang w:val="en-US"/>' +
            '</w:rPr>' +
            '<w:t>p</w:t>' +
          '</w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
        '</w:p>'

$target1.Range.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: the "Форма 116" section's descriptive paragraph ("Тип и количество
# полей совпадают с другими формами.") is replaced by a 2x2 field-description
# table (matching the style already used by the other "Форма N" sections), with
# one field "MAILRANK" (highlighted, description left blank).
# ---------------------------------------------------------------------------
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Тип и количество полей совпадают с другими формами*") {
        $target2 = $p
        break
    }
}

if ($target2 -eq $null) {
    throw "Could not find paragraph containing 'Тип и количество полей совпадают с другими формами'"
}

$xml2 = '<w:tbl xmlns:w="' + $wNs + '">' +
          '<w:tblPr>' +
            '<w:tblStyle w:val="a3"/>' +
            '<w:tblW w:w="0" w:type="auto"/>' +
            '<w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>' +
          '</w:tblPr>' +
          '<w:tblGrid>' +
            '<w:gridCol w:w="2138"/>' +
            '<w:gridCol w:w="2401"/>' +
          '</w:tblGrid>' +
          '<w:tr>' +
            '<w:tc>' +
              '<w:tcPr>' +
                '<w:tcW w:w="2138" w:type="dxa"/>' +
                '<w:shd w:val="clear" w:color="auto" w:fill="EEECE1" w:themeFill="background2"/>' +
              '</w:tcPr>' +
              '<w:p><w:r><w:t>Название поля</w:t></w:r></w:p>' +
            '</w:tc>' +
            '<w:tc>' +
              '<w:tcPr>' +
                '<w:tcW w:w="2401" w:type="dxa"/>' +
                '<w:shd w:val="clear" w:color="auto" w:fill="EEECE1" w:themeFill="background2"/>' +
              '</w:tcPr>' +
              '<w:p><w:r><w:t>Описание</w:t></w:r></w:p>' +
            '</w:tc>' +
          '</w:tr>' +
          '<w:tr>' +
            '<w:tc>' +
              '<w:tcPr>' +
                '<w:tcW w:w="2138" w:type="dxa"/>' +
              '</w:tcPr>' +
              '<w:p>' +
                '<w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr>' +
                '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>MAILRANK</w:t></w:r>' +
              '</w:p>' +
            '</w:tc>' +
            '<w:tc>' +
              '<w:tcPr>' +
                '<w:tcW w:w="2401" w:type="dxa"/>' +
              '</w:tcPr>' +
              '<w:p>' +
                '<w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr>' +
              '</w:p>' +
            '</w:tc>' +
          '</w:tr>' +
        '</w:tbl>'

$target2.Range.InsertXML($xml2) | Out-Null

Write-Output "Edits applied."
